$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset dropped the "Inflammatory-Mac" target cluster entirely, and the TPM-derived
# expression figures were recomputed. Rows 14-17 (the old Resolving-Mac blocks duplicate
# tail plus the now-removed Inflammatory-Mac column) collapse away, leaving a clean 12-row
# (4 sending clusters x 3 target clusters) table in rows 2-13.

# Drop the now-obsolete trailing rows (14-17) first so the sheet ends at row 13.
$ws.Rows.Item(14).Resize(4).Delete() | Out-Null

$data = New-Object 'object[,]' 12,20

$data[0,0] = "ECs"
$data[0,1] = "Epo"
$data[0,2] = "Ephb4"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.07381966666666666
$data[0,7] = 0.221459
$data[0,8] = 0.1284640970637474
$data[0,9] = 0.1284640970637474
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 49.146613
$data[0,13] = 147.439839
$data[0,14] = 0.8214902885327882
$data[0,15] = 0.8214902885327882
$data[0,16] = 3.627986589455666
$data[0,17] = 32.651879305101
$data[0,18] = 0.1055320081630019
$data[0,19] = 0.1055320081630019

$data[1,0] = "ECs"
$data[1,1] = "Epo"
$data[1,2] = "Ephb4"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.07381966666666666
$data[1,7] = 0.221459
$data[1,8] = 0.1284640970637474
$data[1,9] = 0.1284640970637474
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 8.144916333333333
$data[1,13] = 24.434749
$data[1,14] = 0.1361430475126621
$data[1,15] = 0.1361430475126621
$data[1,16] = 0.6012550087545554
$data[1,17] = 5.411295078790999
$data[1,18] = 0.017489493670221
$data[1,19] = 0.017489493670221

$data[2,0] = "ECs"
$data[2,1] = "Epo"
$data[2,2] = "Ephb4"
$data[2,3] = "MuSCs"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.07381966666666666
$data[2,7] = 0.221459
$data[2,8] = 0.1284640970637474
$data[2,9] = 0.1284640970637474
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 2.534635
$data[2,13] = 7.603904999999999
$data[2,14] = 0.04236666395454969
$data[2,15] = 0.04236666395454969
$data[2,16] = 0.1871059108216666
$data[2,17] = 1.683953197395
$data[2,18] = 0.005442595230524438
$data[2,19] = 0.005442595230524438

$data[3,0] = "FAPs"
$data[3,1] = "Epo"
$data[3,2] = "Ephb4"
$data[3,3] = "ECs"
$data[3,4] = 2
$data[3,5] = 0.6666666666666666
$data[3,6] = 0.1698756666666667
$data[3,7] = 0.509627
$data[3,8] = 0.2956247991470493
$data[3,9] = 0.2956247991470493
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 49.146613
$data[3,13] = 147.439839
$data[3,14] = 0.8214902885327882
$data[3,15] = 0.8214902885327882
$data[3,16] = 8.348813647783668
$data[3,17] = 75.139322830053
$data[3,18] = 0.2428529015487571
$data[3,19] = 0.2428529015487571

$data[4,0] = "FAPs"
$data[4,1] = "Epo"
$data[4,2] = "Ephb4"
$data[4,3] = "FAPs"
$data[4,4] = 2
$data[4,5] = 0.6666666666666666
$data[4,6] = 0.1698756666666667
$data[4,7] = 0.509627
$data[4,8] = 0.2956247991470493
$data[4,9] = 0.2956247991470493
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 8.144916333333333
$data[4,13] = 24.434749
$data[4,14] = 0.1361430475126621
$data[4,15] = 0.1361430475126621
$data[4,16] = 1.383623092069222
$data[4,17] = 12.452607828623
$data[4,18] = 0.04024726107619793
$data[4,19] = 0.04024726107619793

$data[5,0] = "FAPs"
$data[5,1] = "Epo"
$data[5,2] = "Ephb4"
$data[5,3] = "MuSCs"
$data[5,4] = 2
$data[5,5] = 0.6666666666666666
$data[5,6] = 0.1698756666666667
$data[5,7] = 0.509627
$data[5,8] = 0.2956247991470493
$data[5,9] = 0.2956247991470493
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 2.534635
$data[5,13] = 7.603904999999999
$data[5,14] = 0.04236666395454969
$data[5,15] = 0.04236666395454969
$data[5,16] = 0.4305728103816667
$data[5,17] = 3.875155293435
$data[5,18] = 0.01252463652209428
$data[5,19] = 0.01252463652209428

$data[6,0] = "MuSCs"
$data[6,1] = "Epo"
$data[6,2] = "Ephb4"
$data[6,3] = "ECs"
$data[6,4] = 2
$data[6,5] = 0.6666666666666666
$data[6,6] = 0.07389766666666667
$data[6,7] = 0.221693
$data[6,8] = 0.1285998359531712
$data[6,9] = 0.1285998359531712
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 49.146613
$data[6,13] = 147.439839
$data[6,14] = 0.8214902885327882
$data[6,15] = 0.8214902885327882
$data[6,16] = 3.631820025269667
$data[6,17] = 32.686380227427
$data[6,18] = 0.1056435163424399
$data[6,19] = 0.1056435163424399

$data[7,0] = "MuSCs"
$data[7,1] = "Epo"
$data[7,2] = "Ephb4"
$data[7,3] = "FAPs"
$data[7,4] = 2
$data[7,5] = 0.6666666666666666
$data[7,6] = 0.07389766666666667
$data[7,7] = 0.221693
$data[7,8] = 0.1285998359531712
$data[7,9] = 0.1285998359531712
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 8.144916333333333
$data[7,13] = 24.434749
$data[7,14] = 0.1361430475126621
$data[7,15] = 0.1361430475126621
$data[7,16] = 0.6018903122285555
$data[7,17] = 5.417012810057
$data[7,18] = 0.01750797357629314
$data[7,19] = 0.01750797357629314

$data[8,0] = "MuSCs"
$data[8,1] = "Epo"
$data[8,2] = "Ephb4"
$data[8,3] = "MuSCs"
$data[8,4] = 2
$data[8,5] = 0.6666666666666666
$data[8,6] = 0.07389766666666667
$data[8,7] = 0.221693
$data[8,8] = 0.1285998359531712
$data[8,9] = 0.1285998359531712
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 2.534635
$data[8,13] = 7.603904999999999
$data[8,14] = 0.04236666395454969
$data[8,15] = 0.04236666395454969
$data[8,16] = 0.1873036123516666
$data[8,17] = 1.685732511165
$data[8,18] = 0.005448346034438223
$data[8,19] = 0.005448346034438223

$data[9,0] = "Resolving-Mac"
$data[9,1] = "Epo"
$data[9,2] = "Ephb4"
$data[9,3] = "ECs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.2570396666666667
$data[9,7] = 0.771119
$data[9,8] = 0.447311267836032
$data[9,9] = 0.447311267836032
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 49.146613
$data[9,13] = 147.439839
$data[9,14] = 0.8214902885327882
$data[9,15] = 0.8214902885327882
$data[9,16] = 12.63262902331567
$data[9,17] = 113.693661209841
$data[9,18] = 0.3674618624785893
$data[9,19] = 0.3674618624785893

$data[10,0] = "Resolving-Mac"
$data[10,1] = "Epo"
$data[10,2] = "Ephb4"
$data[10,3] = "FAPs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 0.2570396666666667
$data[10,7] = 0.771119
$data[10,8] = 0.447311267836032
$data[10,9] = 0.447311267836032
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 8.144916333333333
$data[10,13] = 24.434749
$data[10,14] = 0.1361430475126621
$data[10,15] = 0.1361430475126621
$data[10,16] = 2.093566579347889
$data[10,17] = 18.842099214131
$data[10,18] = 0.06089831918995004
$data[10,19] = 0.06089831918995004

$data[11,0] = "Resolving-Mac"
$data[11,1] = "Epo"
$data[11,2] = "Ephb4"
$data[11,3] = "MuSCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 0.2570396666666667
$data[11,7] = 0.771119
$data[11,8] = 0.447311267836032
$data[11,9] = 0.447311267836032
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 2.534635
$data[11,13] = 7.603904999999999
$data[11,14] = 0.04236666395454969
$data[11,15] = 0.04236666395454969
$data[11,16] = 0.6515017355216666
$data[11,17] = 5.863515619695
$data[11,18] = 0.01895108616749274
$data[11,19] = 0.01895108616749274

$ws.Range("A2:T13").Value = $data
